$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5666.6665
$ws.Range("J69").Value = 5666.6665
$ws.Range("L69").Value = 16999.9995
$ws.Range("N69").Value = -18747.9995

$ws.Range("H72").Value = 5666.6665
$ws.Range("J72").Value = 5666.6665
$ws.Range("L72").Value = 50999.9985
$ws.Range("N72").Value = -59735.9985

$ws.Range("H88").Value = 2501
$ws.Range("I88").Value = 1251.5
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 1251.5
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -845.5
$ws.Range("N88").Value = -5812

$ws.Range("H91").Value = 2501
$ws.Range("I91").Value = 1251.5
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 1251.5
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = 152.5
$ws.Range("N91").Value = -7808

$ws.Range("H139").Value = 58713.332
$ws.Range("J139").Value = 58713.332
$ws.Range("L139").Value = 58713.332
$ws.Range("N139").Value = -68993.33199999999

$ws.Range("H141").Value = 1036.5483
$ws.Range("I141").Value = 533.96155
$ws.Range("K141").Value = 1601.88465
$ws.Range("M141").Value = 3578.11535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3565.0815
$ws.Range("I32").Value = 3534.4597
$ws.Range("K32").Value = 3534.4597
$ws.Range("M32").Value = -3247.4597

$ws.Range("H88").Value = 2389
$ws.Range("I88").Value = 2571.4285
$ws.Range("J88").Value = 1963.3334
$ws.Range("K88").Value = 2571.4285
$ws.Range("L88").Value = 1963.3334
$ws.Range("M88").Value = -2165.4285
$ws.Range("N88").Value = -2775.3334

$ws.Range("H91").Value = 2389
$ws.Range("I91").Value = 2571.4285
$ws.Range("J91").Value = 1963.3334
$ws.Range("K91").Value = 2571.4285
$ws.Range("L91").Value = 1963.3334
$ws.Range("M91").Value = -1167.4285
$ws.Range("N91").Value = -4771.3334

$ws.Range("H132").Value = 592808.5
$ws.Range("I132").Value = 735069.9
$ws.Range("J132").Value = 78479
$ws.Range("K132").Value = 2205209.7
$ws.Range("L132").Value = 235437
$ws.Range("M132").Value = -2202679.7
$ws.Range("N132").Value = -240497

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2120.3333
$ws.Range("I86").Value = 1921.5807
$ws.Range("J86").Value = 2482.7646
$ws.Range("K86").Value = 1921.5807
$ws.Range("L86").Value = 2482.7646
$ws.Range("M86").Value = -798.5807
$ws.Range("N86").Value = -4728.7646

$ws.Range("H89").Value = 2120.3333
$ws.Range("I89").Value = 1921.5807
$ws.Range("J89").Value = 2482.7646
$ws.Range("K89").Value = 9607.9035
$ws.Range("L89").Value = 12413.823
$ws.Range("M89").Value = -3991.9035
$ws.Range("N89").Value = -23645.823

$ws.Range("H105").Value = 1841.3334
$ws.Range("I105").Value = 1847.2727
$ws.Range("J105").Value = 1825
$ws.Range("K105").Value = 1847.2727
$ws.Range("L105").Value = 1825
$ws.Range("M105").Value = -100.2727
$ws.Range("N105").Value = -5319

$ws.Range("H134").Value = 6288293
$ws.Range("I134").Value = 8748576
$ws.Range("J134").Value = 902.55554
$ws.Range("K134").Value = 26245728
$ws.Range("L134").Value = 2707.66662
$ws.Range("M134").Value = -26243193
$ws.Range("N134").Value = -7777.66662

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 460518.22
$ws.Range("I31").Value = 1264.7916
$ws.Range("J31").Value = 1195323.8
$ws.Range("K31").Value = 1264.7916
$ws.Range("L31").Value = 1195323.8
$ws.Range("M31").Value = -969.7916
$ws.Range("N31").Value = -1195913.8

$ws.Range("H34").Value = 460518.22
$ws.Range("I34").Value = 1264.7916
$ws.Range("J34").Value = 1195323.8
$ws.Range("K34").Value = 1264.7916
$ws.Range("L34").Value = 1195323.8
$ws.Range("M34").Value = -1062.7916
$ws.Range("N34").Value = -1195727.8

$ws.Range("H58").Value = 1683.9125
$ws.Range("I58").Value = 767.05554
$ws.Range("K58").Value = 767.05554
$ws.Range("M58").Value = -564.05554

$ws.Range("H136").Value = 1683.9125
$ws.Range("I136").Value = 767.05554
$ws.Range("K136").Value = 2301.16662
$ws.Range("M136").Value = 248.83338

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H105").Value = 4001499
$ws.Range("J105").Value = 4001499
$ws.Range("L105").Value = 12004497
$ws.Range("N105").Value = -12009739

$ws.Range("H110").Value = 3010
$ws.Range("I110").Value = 3000
$ws.Range("J110").Value = 3030
$ws.Range("K110").Value = 9000
$ws.Range("L110").Value = 9090
$ws.Range("M110").Value = -4910
$ws.Range("N110").Value = -17270

$ws.Range("H113").Value = 458.22784
$ws.Range("I113").Value = 447.47916
$ws.Range("J113").Value = 474.87097
$ws.Range("K113").Value = 1342.43748
$ws.Range("L113").Value = 1424.61291
$ws.Range("M113").Value = 827.5625199999999
$ws.Range("N113").Value = -5764.61291

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 24489.898
$ws.Range("I80").Value = 2233.8708
$ws.Range("J80").Value = 62819.723
$ws.Range("K80").Value = 2233.8708
$ws.Range("L80").Value = 62819.723
$ws.Range("M80").Value = -1235.8708
$ws.Range("N80").Value = -64815.723

$ws.Range("H83").Value = 24489.898
$ws.Range("I83").Value = 2233.8708
$ws.Range("J83").Value = 62819.723
$ws.Range("K83").Value = 11169.354
$ws.Range("L83").Value = 314098.615
$ws.Range("M83").Value = -6177.354000000001
$ws.Range("N83").Value = -324082.615

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1778.4
$ws.Range("I61").Value = 1664.5
$ws.Range("J61").Value = 1908.5714
$ws.Range("K61").Value = 1664.5
$ws.Range("L61").Value = 1908.5714
$ws.Range("M61").Value = -1462.5
$ws.Range("N61").Value = -2312.5714

$ws.Range("H68").Value = 1862.75
$ws.Range("I68").Value = 1660.4
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 1660.4
$ws.Range("L68").Value = 2200
$ws.Range("M68").Value = -911.4000000000001
$ws.Range("N68").Value = -3698

$ws.Range("H71").Value = 1862.75
$ws.Range("I71").Value = 1660.4
$ws.Range("J71").Value = 2200
$ws.Range("K71").Value = 8302
$ws.Range("L71").Value = 11000
$ws.Range("M71").Value = -4558
$ws.Range("N71").Value = -18488

$ws.Range("H113").Value = 1778.4
$ws.Range("I113").Value = 1664.5
$ws.Range("J113").Value = 1908.5714
$ws.Range("K113").Value = 1664.5
$ws.Range("L113").Value = 1908.5714
$ws.Range("M113").Value = 505.5
$ws.Range("N113").Value = -6248.5714

$ws.Range("H122").Value = 2068.889
$ws.Range("I122").Value = 1952.5
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 5857.5
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -3407.5
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 466.85715
$ws.Range("I113").Value = 451
$ws.Range("J113").Value = 488
$ws.Range("K113").Value = 1353
$ws.Range("L113").Value = 1464
$ws.Range("M113").Value = 817
$ws.Range("N113").Value = -5804

$ws.Range("H122").Value = 2637.762
$ws.Range("I122").Value = 1749.5714
$ws.Range("J122").Value = 4414.143
$ws.Range("K122").Value = 5248.7142
$ws.Range("L122").Value = 13242.429
$ws.Range("M122").Value = -2798.7142
$ws.Range("N122").Value = -18142.429

$ws.Range("H126").Value = 1667.3334
$ws.Range("I126").Value = 1432.5454
$ws.Range("J126").Value = 4250
$ws.Range("K126").Value = 4297.6362
$ws.Range("L126").Value = 12750
$ws.Range("M126").Value = -1827.6362
$ws.Range("N126").Value = -17690

$ws.Range("H132").Value = 4319.171
$ws.Range("I132").Value = 4762.6665
$ws.Range("J132").Value = 1126
$ws.Range("K132").Value = 14287.9995
$ws.Range("L132").Value = 3378
$ws.Range("M132").Value = -11757.9995
$ws.Range("N132").Value = -8438
